$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '244.27'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-1.17%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '27.11'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '2.84%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.024'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-1.11%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05664'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.83%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.471'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-0.33%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8220'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '1.17%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8416'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-0.44%'
$ws.Range('B9').Value = 'One'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0005990'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '0.17%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1324'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-1.50%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06916'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-0.62%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.02883'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '2.21%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09391'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.09%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001520'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.23%'
$ws.Range('B15').Value = 'CoinExToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.04134'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-11.44%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006213'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-0.45%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.511'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-1.86%'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-1.88%'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '9.10%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.03151'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-0.40%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.1255'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-4.98%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.580'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-4.75%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001222'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-2.40%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.003870'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-16.30%'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '2.05%'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '-25.77%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03669'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '0.13%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006070'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '77.07%'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-21.97%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002277'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-14.41%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01134'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '30.53%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005314'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '0.45%'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.01%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1015'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-15.41%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002571'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '24.48%'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '0.01%'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.01%'

Write-Host "Applied crypto price/volume updates"
